$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.801.81"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.68%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.452.11"
$ws.Range("D3").Style = "Normal"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "160.12"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.36%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.454.59"
$ws.Range("D8").Style = "Normal"
$ws.Range("E9").Value = "  +9.28%  "
$ws.Range("E10").Value = "  -2.79%  "
$ws.Range("E11").Value = "  +2.55%  "
$ws.Range("E12").Value = "  +0.66%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.050.16"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.61%  "
$ws.Range("E14").Value = "  -2.45%  "
$ws.Range("E15").Value = "  +4.36%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "28.19"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.63%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.832.07"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.59%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.454.10"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.58%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.35"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.31"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.40%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "387.37"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.94%  "
$ws.Range("E22").Value = "  -3.60%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "73.22"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.57%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.545"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.88%  "
$ws.Range("E25").Value = "  +0.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000124"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +15.91%  "
$ws.Range("E27").Value = "  +1.85%  "
$ws.Range("E28").Value = "  -0.49%  "
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.23"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +8.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.43"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.23%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.05"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.71"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.98%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.52"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.46%  "
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("E36").Value = "  +4.58%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "163.53"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.94%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.51"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.77%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.90"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.82%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.016.89"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.51%  "
$ws.Range("E41").Value = "  -2.13%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "27.25"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.86%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.58"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.28%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "42.93"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.58%  "
$ws.Range("E45").Value = "  -1.54%  "
$ws.Range("E46").Value = "  +0.76%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "24.55"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +8.83%  "
$ws.Range("E48").Value = "  +0.26%  "
$ws.Range("E49").Value = "  +6.14%  "
$ws.Range("B50").Value = "dogwifhat"
$ws.Range("C50").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.17"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.31%  "
$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.59"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.31%  "
